$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.090.89'
$ws.Range('E2').Value = '  -4.86%  '
$ws.Range('D3').Value = '3.230.05'
$ws.Range('E3').Value = '  -8.24%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '592.20'
$ws.Range('E5').Value = '  -2.45%  '
$ws.Range('D6').Value = '152.52'
$ws.Range('E6').Value = '  -12.64%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.214.43'
$ws.Range('E8').Value = '  -8.58%  '
$ws.Range('E9').Value = '  -11.45%  '
$ws.Range('E10').Value = '  -12.32%  '
$ws.Range('D11').Value = '6.51'
$ws.Range('E11').Value = '  -10.53%  '
$ws.Range('D12').Value = '0.494'
$ws.Range('E12').Value = '  -16.11%  '
$ws.Range('D13').Value = '38.94'
$ws.Range('E13').Value = '  -16.03%  '
$ws.Range('E14').Value = '  -12.20%  '
$ws.Range('D15').Value = '3.757.25'
$ws.Range('E15').Value = '  -8.18%  '
$ws.Range('D16').Value = '67.212.74'
$ws.Range('E16').Value = '  -4.74%  '
$ws.Range('D17').Value = '3.240.60'
$ws.Range('E17').Value = '  -7.96%  '
$ws.Range('E18').Value = '  -4.49%  '
$ws.Range('D19').Value = '531.76'
$ws.Range('E19').Value = '  -13.24%  '
$ws.Range('D20').Value = '7.11'
$ws.Range('E20').Value = '  -14.61%  '
$ws.Range('D21').Value = '14.89'
$ws.Range('E21').Value = '  -15.14%  '
$ws.Range('E22').Value = '  -14.17%  '
$ws.Range('E23').Value = '  -12.39%  '
$ws.Range('D24').Value = '13.80'
$ws.Range('E24').Value = '  -11.70%  '
$ws.Range('D25').Value = '85.60'
$ws.Range('E25').Value = '  -14.11%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -14.49%  '
$ws.Range('D28').Value = '2.18'
$ws.Range('E28').Value = '  -15.38%  '
$ws.Range('D29').Value = '8.06'
$ws.Range('E29').Value = '  -11.00%  '
$ws.Range('D30').Value = '29.09'
$ws.Range('E30').Value = '  -15.31%  '
$ws.Range('D31').Value = '2.66'
$ws.Range('E31').Value = '  -10.66%  '
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  -9.99%  '
$ws.Range('D33').Value = '543.92'
$ws.Range('E33').Value = '  -15.64%  '
$ws.Range('D34').Value = '5.74'
$ws.Range('E34').Value = '  -16.31%  '
$ws.Range('D35').Value = '6.43'
$ws.Range('E35').Value = '  -20.39%  '
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('D37').Value = '53.53'
$ws.Range('E37').Value = '  -5.85%  '
$ws.Range('E38').Value = '  -11.07%  '
$ws.Range('D39').Value = '0.0860'
$ws.Range('E39').Value = '  -13.86%  '
$ws.Range('E40').Value = '  -14.02%  '
$ws.Range('E41').Value = '  -13.08%  '
$ws.Range('D42').Value = '2.934.93'
$ws.Range('E42').Value = '  -12.92%  '
$ws.Range('D43').Value = '2.67'
$ws.Range('E43').Value = '  -25.25%  '
$ws.Range('E44').Value = '  -15.32%  '
$ws.Range('D45').Value = '0.0₃0585'
$ws.Range('E45').Value = '  -21.87%  '
$ws.Range('D46').Value = '2.43'
$ws.Range('E46').Value = '  -16.09%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '26.53'
$ws.Range('E47').Value = '  -17.79%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').Value = '2.16'
$ws.Range('E48').Value = '  -16.00%  '
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('E50').Value = '  -11.84%  '
$ws.Range('D51').Value = '118.25'
$ws.Range('E51').Value = '  -11.39%  '
